$wb = $excel.ActiveWorkbook
Write-Output $excel.CalculationVersion
